$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2025-12-26 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-12-27 Saturday", 2)

# Update the multiplication problems in the table (single table, 20 rows x 5 cols,
# with content only in rows 1, 5, 10, 15, 20).
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text  = "41×93="
$t.Cell(1, 2).Range.Text  = "79×46="
$t.Cell(1, 3).Range.Text  = "16×95="
$t.Cell(1, 4).Range.Text  = "73×11="
$t.Cell(1, 5).Range.Text  = "74×85="

$t.Cell(5, 1).Range.Text  = "21×42="
$t.Cell(5, 2).Range.Text  = "37×47="
$t.Cell(5, 3).Range.Text  = "55×32="
$t.Cell(5, 4).Range.Text  = "11×27="
$t.Cell(5, 5).Range.Text  = "17×60="

$t.Cell(10, 1).Range.Text = "93×86="
$t.Cell(10, 2).Range.Text = "71×36="
$t.Cell(10, 3).Range.Text = "59×12="
$t.Cell(10, 4).Range.Text = "68×11="
$t.Cell(10, 5).Range.Text = "38×45="

$t.Cell(15, 1).Range.Text = "19×76="
$t.Cell(15, 2).Range.Text = "17×32="
$t.Cell(15, 3).Range.Text = "27×44="
$t.Cell(15, 4).Range.Text = "36×35="
$t.Cell(15, 5).Range.Text = "33×75="

$t.Cell(20, 1).Range.Text = "99×50="
$t.Cell(20, 2).Range.Text = "17×76="
$t.Cell(20, 3).Range.Text = "16×59="
$t.Cell(20, 4).Range.Text = "37×66="
$t.Cell(20, 5).Range.Text = "99×87="
